# TAYD Quarterly Financials update
# Inserts 3 new quarterly columns (new periods ending 2019-02-28, 2018-11-30, 2018-08-31)
# before the existing data (old column D), shifting the prior 8 quarters right by 3 columns,
# and fills in the figures for the 3 new quarters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 blank columns at D:F; existing D:K data (and all its formatting) shifts to G:N.
$ws.Range("D1:F1").EntireColumn.Insert()

# 2) The newly inserted D:F columns come in unformatted; copy number formats (date / number)
#    from column G - the just-shifted former column D - so D:F match the rest of the table.
$ws.Range("G5:G102").Copy()
$ws.Range("D5:F102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) New quarter figures, keyed by row number -> (D, E, F) values
#    (D = period ending 2019-02-28, E = 2018-11-30, F = 2018-08-31)
$newData = @{
    7 = @(43524, 43434, 43343)
    8 = @(7800, 9500, 7300)
    9 = @(5800, 7300, 5000)
    10 = @(2000, 2200, 2300)
    11 = @($null, $null, $null)
    12 = @("NA", "NA", "NA")
    13 = @(0, 0, 0)
    14 = @(0, 0, 0)
    15 = @(0, 0, 0)
    16 = @($null, $null, $null)
    17 = @(7300, 9000, 6400)
    18 = @(500, 500, 900)
    19 = @($null, $null, $null)
    20 = @(0, 0, 0)
    21 = @(800, 700, 1200)
    22 = @(0, 0, 0)
    23 = @(500, 500, 900)
    24 = @(100, 100, 200)
    25 = @(0, 0, 0)
    26 = @(400, 400, 700)
    27 = @(400, 400, 700)
    28 = @(0, 0, 0)
    29 = @("NA", "NA", "NA")
    30 = @(0, 0, 0)
    31 = @(0, 0, 0)
    32 = @(0, 0, 0)
    33 = @(400, 400, 700)
    34 = @(0, 0, 0)
    35 = @(400, 400, 700)
    38 = @(43524, 43434, 43343)
    39 = @($null, $null, $null)
    40 = @($null, $null, $null)
    41 = @(4300, 5000, 7100)
    42 = @(1100, 1000, 1000)
    43 = @(4800, 5000, 4800)
    44 = @(18700, 18000, 17000)
    45 = @(700, 500, 300)
    46 = @(29600, 29500, 30300)
    47 = @(0, 0, 0)
    48 = @(9600, 9700, 9700)
    49 = @(0, 0, 0)
    50 = @(0, 0, 0)
    51 = @(0, 0, 0)
    52 = @(1100, 1300, 1300)
    53 = @(0, 0, 0)
    54 = @(40200, 40500, 41400)
    55 = @($null, $null, $null)
    56 = @($null, $null, $null)
    57 = @(1600, 2100, 1900)
    58 = @(0, 0, 0)
    59 = @(3400, 3600, 5100)
    60 = @(5000, 5700, 7000)
    61 = @(0, 0, 0)
    62 = @(0, 0, 0)
    63 = @(0, 0, 0)
    64 = @(0, 0, 0)
    65 = @(0, 0, 0)
    66 = @(5000, 5700, 7000)
    67 = @($null, $null, $null)
    68 = @(0, 0, 0)
    69 = @(0, 0, 0)
    70 = @(0, 0, 0)
    71 = @(0, 0, 0)
    72 = @(28500, 28100, 27700)
    73 = @(0, 0, 0)
    74 = @(0, 0, 0)
    75 = @(0, 0, 0)
    76 = @(35300, 34800, 34400)
    77 = @(0, 0, 0)
    79 = @($null, $null, $null)
    80 = @(43524, 43434, 43343)
    81 = @(400, 400, 700)
    82 = @($null, $null, $null)
    83 = @(300, 300, 300)
    84 = @(0, 0, 0)
    85 = @(0, 0, 0)
    86 = @(0, 0, 0)
    87 = @(0, 0, 0)
    88 = @(0, 0, 0)
    89 = @(-500, -1900, 4400)
    90 = @($null, $null, $null)
    91 = @(-100, -200, -100)
    92 = @(0, 0, 0)
    93 = @(0, 0, 0)
    94 = @(-100, -200, -100)
    95 = @($null, $null, $null)
    96 = @(0, 0, 0)
    97 = @(0, 0, 0)
    98 = @(0, 0, 0)
    99 = @(0, 0, 0)
    100 = @(0, 0, 0)
    101 = @(0, 0, 0)
    102 = @(-700, -2200, 4300)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    if ($vals[0] -ne $null) { $ws.Range("D$row").Value = $vals[0] }
    if ($vals[1] -ne $null) { $ws.Range("E$row").Value = $vals[1] }
    if ($vals[2] -ne $null) { $ws.Range("F$row").Value = $vals[2] }
}

# 4) Row 62 ("Other Liabilities"): the four oldest-of-the-shifted quarters (now G:J, formerly
#    D:G) were corrected from "NA" to 0 in this revision; the oldest one (now K, formerly H)
#    keeps its "NA".
$ws.Range("G62:J62").Value = 0
